# Applies the commit's textual + structural changes to the purchase-order
# style document:
#   1. Two stray placeholder strings corrected.
#   2. The document date corrected.
#   3. Two of the three (duplicate) line-item rows removed from the
#      products table, and the remaining item's description/price updated.
#   4. The Gross/Sub/Grand total cells updated to match the new total.

$d = $word.ActiveDocument

# 1) Stray placeholder text fixes -------------------------------------------------
$d.Content.Find.Execute("sdadsa", $true, $false, $false, $false, $false, `
    $true, 1, $false, "asdasd", 2) | Out-Null

$d.Content.Find.Execute("sadasd", $true, $false, $false, $false, $false, `
    $true, 1, $false, "asdsadsa", 2) | Out-Null

# 2) Document date ----------------------------------------------------------------
$d.Content.Find.Execute("28 มี.ค. 2023", $true, $false, $false, $false, $false, `
    $true, 1, $false, "3 ก.พ. 2023 ", 2) | Out-Null

# 3) Remove the duplicate line-item rows (items "2" and "3") from the
#    products table, keeping only item "1" (which gets re-described below).
#    The products table is the 4th table in the document.
$productsTable = $d.Tables.Item(4)
$productsTable.Rows.Item(5).Delete()
$productsTable.Rows.Item(4).Delete()

# Update the surviving line item's description and price.
$d.Content.Find.Execute("ชุดหน้าใส", $true, $false, $false, $false, $false, `
    $true, 1, $false, "โลชั่น เรตินอล", 2) | Out-Null

$d.Content.Find.Execute("4100.0", $true, $false, $false, $false, $false, `
    $true, 1, $false, "1190.0", 2) | Out-Null

# 4) Update the Gross Total / Sub Total / Grand Total cells (all three still
#    read the old combined total) to the new single-item total.
while ($d.Content.Find.Execute("12300.0", $true, $false, $false, $false, $false, `
    $true, 1, $false, "1190.0", 2)) { }
